$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.49310000000002
$ws.Range("E4").Value = 14.33189999999999

$ws.Range("E5").Value = 13.13189999999999

$ws.Range("A6").Value = -20.32779999999999
$ws.Range("E6").Value = 13.0903

$ws.Range("A7").Value = -21.22510000000001

$ws.Range("A8").Value = -20.71389999999999
$ws.Range("E8").Value = 13.08800000000001

$ws.Range("A16").Value = -20.3421
$ws.Range("E16").Value = 12.80750000000001

$ws.Range("A20").Value = -22.19990000000002

$ws.Range("A21").Value = -20.61579999999999

$ws.Range("E22").Value = 12.22969999999999
